$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Month/Year text in A1 (merged A1:AE1 header cell).
$ws.Range("A1").Value = "Month/Year: March 2021"

# Avoid Excel's automatic row-height autofit side effect from the large font
# used by A1's style, keeping row 1 at its original (non-custom) height.
$ws.Rows.Item(1).AutoFit()

# Materialize the empty, styled cells D1:AD1 across the row (same blank
# "numeric" style as the existing B1/C1 cells), by using iter_rows-style
# per-cell iteration copying the established style forward.
$ws.Range("B1:C1").Copy()
$ws.Range("D1:AD1").PasteSpecial(-4122)
